$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 5 de Septiembre de 2020 a las 06:31"

# Row 6 - India
$ws.Range("B6").Value = 4023179
$ws.Range("C6").Value = 2940
$ws.Range("D6").Value = 3107223
$ws.Range("E6").Value = 846321

# Row 20 - Pakistan
$ws.Range("B20").Value = 298025
$ws.Range("C20").Value = 513
$ws.Range("D20").Value = 282553
$ws.Range("E20").Value = 9132
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 6340

# Row 124 - Tailandia
$ws.Range("B124").Value = 3438
$ws.Range("C124").Value = 7
$ws.Range("D124").Value = 3279
$ws.Range("E124").Value = 101

# Row 187 - Butan
$ws.Range("B187").Value = 228
$ws.Range("C187").Value = 1
$ws.Range("E187").Value = 78
